$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 106, shifting existing rows
# 106-113 down to 107-114 (preserving their values/formatting).
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with the weekly Cilantro price entry.
$ws.Range("A106").Value = 1
$ws.Range("B106").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C106").Value = "Arica y Parinacota"
$ws.Range("D106").Value = 45041
$ws.Range("D106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E106").Value = 15
$ws.Range("F106").Value = 100112040
$ws.Range("G106").Value = "Cilantro"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 300
$ws.Range("K106").Value = 1500
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = 1750
$ws.Range("N106").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O106").Value = "Región de Arica y Parinacota"
$ws.Range("P106").Value = 875
$ws.Range("Q106").Value = 2
$ws.Range("R106").Value = "Hortaliza"
